$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 - Key Stage 4 (KS4) destinations: roll the "next"/"latest" period data
# forward one year (now 21/22 learners) and point the source link at the new
# permalink.
$ws.Range("B11").Value = "<a href = 'https://explore-education-statistics.service.gov.uk/data-tables/permalink/f5995891-4e9a-439d-0d58-08dbd1400c33'>National Pupil Database</a>"
$ws.Range("C11").Value = "Aug 2021 -  Jul 2022 (20/21 learners) (19/10/23)"
$ws.Range("D11").Value = "Aug 2022 -  Jul 2023 (21/22 learners) (Oct 24)"

# Row 12 - Key Stage 5 (KS5) destinations: same period roll-forward, different
# source permalink.
$ws.Range("B12").Value = "<a href = 'https://explore-education-statistics.service.gov.uk/data-tables/permalink/b0424f32-b140-44a8-a039-08dbd1466a44'>National Pupil Database</a>"
$ws.Range("C12").Value = "Aug 2021 -  Jul 2022 (20/21 learners) (19/10/23)"
$ws.Range("D12").Value = "Aug 2022 -  Jul 2023 (21/22 learners) (Oct 24)"

# Reflect the author's last active cell selection on the sheet.
$ws.Range("B12").Select()
